# Auto-generated Excel COM-interop script
# Updates Price (D) and Volume(1h) (E) columns per the Jan 12 2023 symbol-list refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'285.46"
$ws.Range("E2").Value = "'3.07%"
$ws.Range("D3").Value = "'28.62"
$ws.Range("E3").Value = "'4.51%"
$ws.Range("D4").Value = "'5.032"
$ws.Range("E4").Value = "'3.25%"
$ws.Range("D5").Value = "'0.06490"
$ws.Range("E5").Value = "'1.15%"
$ws.Range("D6").Value = "'7.231"
$ws.Range("E6").Value = "'4.19%"
$ws.Range("D7").Value = "'1.339"
$ws.Range("E7").Value = "'13.40%"
$ws.Range("D8").Value = "'0.9121"
$ws.Range("E8").Value = "'4.20%"
$ws.Range("D9").Value = "'0.1540"
$ws.Range("E9").Value = "'0.34%"
$ws.Range("D10").Value = "'0.06487"
$ws.Range("E10").Value = "'25.70%"
$ws.Range("D11").Value = "'0.07595"
$ws.Range("E11").Value = "'2.22%"
$ws.Range("D12").Value = "'0.02981"
$ws.Range("E12").Value = "'1.24%"
$ws.Range("D13").Value = "'0.08955"
$ws.Range("E13").Value = "'-0.30%"
$ws.Range("D14").Value = "'0.001602"
$ws.Range("E14").Value = "'2.31%"
$ws.Range("D15").Value = "'0.0006519"
$ws.Range("E15").Value = "'2.56%"
$ws.Range("D16").Value = "'0.006034"
$ws.Range("E16").Value = "'-1.38%"
$ws.Range("D17").Value = "'3.460"
$ws.Range("E17").Value = "'-0.55%"
$ws.Range("E18").Value = "'1.79%"
$ws.Range("D20").Value = "'0.3178"
$ws.Range("E20").Value = "'1.35%"
$ws.Range("D21").Value = "'0.1342"
$ws.Range("E21").Value = "'1.50%"
$ws.Range("D22").Value = "'3.974"
$ws.Range("E22").Value = "'1.83%"
$ws.Range("D24").Value = "'0.04458"
$ws.Range("E24").Value = "'1.12%"
$ws.Range("E25").Value = "'0.56%"
$ws.Range("D26").Value = "'0.004329"
$ws.Range("E28").Value = "'-9.18%"
$ws.Range("D29").Value = "'0.0001635"
$ws.Range("E29").Value = "'-15.74%"
$ws.Range("D40").Value = "'0.04143"
$ws.Range("E40").Value = "'-0.54%"
$ws.Range("D41").Value = "'0.006777"
$ws.Range("E41").Value = "'-0.45%"
$ws.Range("D42").Value = "'0.1233"
$ws.Range("E42").Value = "'5.15%"
$ws.Range("D43").Value = "'0.002121"
$ws.Range("E43").Value = "'4.03%"
$ws.Range("D44").Value = "'0.01188"
$ws.Range("E44").Value = "'3.81%"
$ws.Range("D45").Value = "'0.00005383"
$ws.Range("E45").Value = "'1.33%"
$ws.Range("E46").Value = "'-0.02%"
$ws.Range("D47").Value = "'1.933"
$ws.Range("E47").Value = "'14.73%"
